$wb = $excel.ActiveWorkbook

# "Topic-contributed Sessions" sheet: update the organizer name for the
# "Generalized pairwise comparisons" session and make this sheet the active tab.
$ws2 = $wb.Worksheets.Item("Topic-contributed Sessions")
$ws2.Range("B9").Value = "Arne Bathke, Johan Verbeeck"

$ws2.Activate()
$ws2.Range("B9").Select()
